$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Tree 1" sheet: the bracket-drawing grid (columns B:I) shifts down by one
#    row relative to column A (rows 2-15 -> rows 3-16), while the pool-name
#    labels drawn in that grid get re-assigned to new pools.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Tree 1")

# Each entry: source cell (supplies formatting), destination cell, content
# type ("s" = shared string literal, "n" = number, $null = blank/format-only),
# and the literal content to place at the destination.
$moves = @(
    @("E2","E3","s","Pool A.1"),
    @("F3","F4",$null,$null),
    @("G3","G4",$null,$null),
    @("G4","G5","n",3),
    @("C5","C6","s","Pool C.1"),
    @("G5","G6",$null,$null),
    @("H5","H6",$null,$null),
    @("I5","I6",$null,$null),
    @("D6","D7",$null,$null),
    @("E6","E7","n",1),
    @("F6","F7",$null,$null),
    @("G6","G7",$null,$null),
    @("I6","I7",$null,$null),
    @("C7","C8","s","Pool B.2"),
    @("D7","D8",$null,$null),
    @("E7","E8",$null,$null),
    @("I7","I8",$null,$null),
    @("I8","I9","n",5),
    @("I9","I10",$null,$null),
    @("E10","E11","s","Pool B.1"),
    @("I10","I11",$null,$null),
    @("F11","F12",$null,$null),
    @("G11","G12",$null,$null),
    @("I11","I12",$null,$null),
    @("G12","G13","n",4),
    @("H12","H13",$null,$null),
    @("I12","I13",$null,$null),
    @("C13","C14","s","Pool A.2"),
    @("G13","G14",$null,$null),
    @("D14","D15",$null,$null),
    @("E14","E15","n",2),
    @("F14","F15",$null,$null),
    @("G14","G15",$null,$null),
    @("C15","C16","s","Pool C.2"),
    @("D15","D16",$null,$null),
    @("E15","E16",$null,$null)
)

# Phase 1: snapshot the formatting of every source cell by copying it onto
# the destination cell (format only) BEFORE any cell is cleared or
# overwritten, since several destinations double as sources for the next row.
foreach ($m in $moves) {
    $srcRef = $m[0]
    $dstRef = $m[1]
    $ws.Range($srcRef).Copy()
    $ws.Range($dstRef).PasteSpecial(-4122)
}

# Phase 2: cells that only ever acted as a source (never a destination) are
# fully vacated by the shift - clear content AND formatting so no stray cell
# is left behind.
$dstRefs = @{}
foreach ($m in $moves) {
    $dstRefs[$m[1]] = $true
}
$srcRefsSeen = @{}
foreach ($m in $moves) {
    $ref = $m[0]
    if ($srcRefsSeen.ContainsKey($ref)) { continue }
    $srcRefsSeen[$ref] = $true
    if (-not $dstRefs.ContainsKey($ref)) {
        $ws.Range($ref).Clear()
    }
}

# Phase 3: every destination cell's stale content (if any) is cleared, then
# the final literal content is written in.
foreach ($ref in $dstRefs.Keys) {
    $ws.Range($ref).ClearContents()
}
foreach ($m in $moves) {
    $dstRef = $m[1]
    $kind = $m[2]
    $content = $m[3]
    if ($kind -eq "s") {
        $ws.Range($dstRef).Value = $content
    } elseif ($kind -eq "n") {
        $ws.Range($dstRef).Value = $content
    }
}

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2) "Time Estimator" sheet: "Matches per pool" (C2) goes from 3 to 6.
# ---------------------------------------------------------------------------
$wsTime = $wb.Worksheets.Item("Time Estimator")
$wsTime.Range("C2").Value = 6

# ---------------------------------------------------------------------------
# 3) "Names to Print" sheet: the pool-slot label formulas are re-assigned to
#    different pools/cells in rows 5, 12, 28 and 35.
# ---------------------------------------------------------------------------
$wsNames = $wb.Worksheets.Item("Names to Print")

$wsNames.Range("A5").Formula = "=CONCATENATE(""Pool C.1 "",'Pool Matches'!G125)"
$wsNames.Range("G5").Formula = "=CONCATENATE(""Pool B.2 "",'Pool Matches'!O43)"
$wsNames.Range("I5").Formula = "=CONCATENATE(""Pool A.2 "",'Pool Matches'!G82)"
$wsNames.Range("O5").Formula = "=CONCATENATE(""Pool C.2 "",'Pool Matches'!G126)"

$wsNames.Range("A12").Formula = "=CONCATENATE(""Pool C.1 "",'Pool Matches'!G125)"
$wsNames.Range("G12").Formula = "=CONCATENATE(""Pool B.2 "",'Pool Matches'!O43)"
$wsNames.Range("I12").Formula = "=CONCATENATE(""Pool A.2 "",'Pool Matches'!G82)"
$wsNames.Range("O12").Formula = "=CONCATENATE(""Pool C.2 "",'Pool Matches'!G126)"

$wsNames.Range("I28").Formula = "=CONCATENATE(""Pool B.1 "",'Pool Matches'!O42)"
$wsNames.Range("I35").Formula = "=CONCATENATE(""Pool B.1 "",'Pool Matches'!O42)"
